$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value (text, to preserve the original inline-string/text cell type)
$updates = [ordered]@{
    'D2' = '309.73'
    'E2' = '-2.46%'
    'D3' = '52.62'
    'E3' = '9.07%'
    'D4' = '5.108'
    'E4' = '-2.92%'
    'D5' = '0.07784'
    'E5' = '-2.53%'
    'D6' = '4.492'
    'E6' = '-2.20%'
    'D7' = '1.359'
    'E7' = '-4.98%'
    'D8' = '1.580'
    'E8' = '-3.98%'
    'D9' = '0.1228'
    'E9' = '-4.10%'
    'D10' = '0.2000'
    'E10' = '3.33%'
    'D11' = '0.09511'
    'E11' = '2.53%'
    'D12' = '0.04722'
    'E12' = '2.41%'
    'D13' = '0.1045'
    'E13' = '0.09%'
    'D14' = '0.001261'
    'E14' = '-4.73%'
    'D15' = '0.005800'
    'E15' = '-1.82%'
    'E16' = '2,009.50%'
    'E17' = '0.13%'
    'E18' = '-0.74%'
    'E19' = '0.06%'
    'D20' = '8.008'
    'E20' = '-1.39%'
    'E21' = '-2.90%'
    'E22' = '-0.41%'
    'D23' = '0.04180'
    'E23' = '0.20%'
    'E24' = '-3.79%'
    'D25' = '0.003946'
    'E25' = '-6.74%'
    'E26' = '-0.02%'
    'E38' = '-3.01%'
    'D39' = '0.05889'
    'E39' = '4.40%'
    'D40' = '0.01102'
    'E40' = '34.09%'
    'D41' = '0.008090'
    'E41' = '-0.70%'
    'E42' = '-0.85%'
    'D43' = '0.008240'
    'E43' = '7.12%'
    'D44' = '0.008440'
    'E44' = '-0.67%'
    'D45' = '0.3114'
    'E45' = '-10.23%'
    'D46' = '0.00007349'
    'E46' = '6.22%'
    'E47' = '-0.10%'
    'D48' = '0.05670'
    'E48' = '3.18%'
    'D49' = '0.002622'
    'E49' = '-34.63%'
    'D50' = '0.00002102'
    'E50' = '-0.10%'
    'E51' = '-0.10%'
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    # Force text format first so Excel keeps these numeric-looking / percent-looking
    # strings as text instead of auto-converting them to Number cells.
    $range.NumberFormat = "@"
    $range.Value = $updates[$cell]
}
